$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 524.02856
$ws.Range("J17").Value = 524.02856
$ws.Range("L17").Value = 1572.08568
$ws.Range("N17").Value = -1908.08568
$ws.Range("H40").Value = 3987.4443
$ws.Range("I40").Value = 4173
$ws.Range("J40").Value = 3894.6667
$ws.Range("K40").Value = 4173
$ws.Range("L40").Value = 3894.6667
$ws.Range("M40").Value = -3998
$ws.Range("N40").Value = -4244.6667
$ws.Range("H69").Value = 8534
$ws.Range("I69").Value = 6951.8
$ws.Range("J69").Value = 12489.5
$ws.Range("K69").Value = 20855.4
$ws.Range("L69").Value = 37468.5
$ws.Range("M69").Value = -19981.4
$ws.Range("N69").Value = -39216.5
$ws.Range("H72").Value = 8534
$ws.Range("I72").Value = 6951.8
$ws.Range("J72").Value = 12489.5
$ws.Range("K72").Value = 62566.2
$ws.Range("L72").Value = 112405.5
$ws.Range("M72").Value = -58198.2
$ws.Range("N72").Value = -121141.5
$ws.Range("H86").Value = 1347.6
$ws.Range("I86").Value = 1213.25
$ws.Range("K86").Value = 1213.25
$ws.Range("M86").Value = -90.25
$ws.Range("H88").Value = 721421.0600000001
$ws.Range("I88").Value = 6686.75
$ws.Range("K88").Value = 6686.75
$ws.Range("M88").Value = -6280.75
$ws.Range("H89").Value = 1347.6
$ws.Range("I89").Value = 1213.25
$ws.Range("K89").Value = 6066.25
$ws.Range("M89").Value = -450.25
$ws.Range("H91").Value = 721421.0600000001
$ws.Range("I91").Value = 6686.75
$ws.Range("K91").Value = 6686.75
$ws.Range("M91").Value = -5282.75
$ws.Range("H106").Value = 61708
$ws.Range("I106").Value = 93347
$ws.Range("K106").Value = 93347
$ws.Range("M106").Value = -92716
$ws.Range("H112").Value = 1455.5
$ws.Range("J112").Value = 1622.7241
$ws.Range("L112").Value = 4868.1723
$ws.Range("N112").Value = -7084.1723
$ws.Range("H137").Value = 1435.8096
$ws.Range("I137").Value = 1603.6897
$ws.Range("J137").Value = 1061.3077
$ws.Range("K137").Value = 4811.0691
$ws.Range("L137").Value = 3183.9231
$ws.Range("M137").Value = -2261.0691
$ws.Range("N137").Value = -8283.9231
$ws.Range("H138").Value = 2077.342
$ws.Range("J138").Value = 2416.6553
$ws.Range("L138").Value = 7249.965899999999
$ws.Range("N138").Value = -17529.9659
$ws.Range("H139").Value = 99979.5
$ws.Range("J139").Value = 99979.5
$ws.Range("L139").Value = 99979.5
$ws.Range("N139").Value = -110259.5
$ws.Range("H140").Value = 86665
$ws.Range("J140").Value = 99995
$ws.Range("L140").Value = 99995
$ws.Range("N140").Value = -110355
$ws.Range("H141").Value = 3108.2222
$ws.Range("I141").Value = 3153.9614
$ws.Range("K141").Value = 9461.8842
$ws.Range("M141").Value = -4281.8842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2971.1724
$ws.Range("I32").Value = 1320.0625
$ws.Range("K32").Value = 1320.0625
$ws.Range("M32").Value = -1033.0625
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H45").Value = 1686.6897
$ws.Range("I45").Value = 1494.3125
$ws.Range("J45").Value = 1923.4615
$ws.Range("K45").Value = 1494.3125
$ws.Range("L45").Value = 1923.4615
$ws.Range("M45").Value = -1117.3125
$ws.Range("N45").Value = -2677.4615
$ws.Range("H74").Value = 4234.3184
$ws.Range("I74").Value = 4032.5
$ws.Range("K74").Value = 4032.5
$ws.Range("M74").Value = -3158.5
$ws.Range("H77").Value = 4234.3184
$ws.Range("I77").Value = 4032.5
$ws.Range("K77").Value = 20162.5
$ws.Range("M77").Value = -15794.5
$ws.Range("H88").Value = 2199.2
$ws.Range("I88").Value = 1912.8889
$ws.Range("J88").Value = 2360.25
$ws.Range("K88").Value = 1912.8889
$ws.Range("L88").Value = 2360.25
$ws.Range("M88").Value = -1506.8889
$ws.Range("N88").Value = -3172.25
$ws.Range("H91").Value = 2199.2
$ws.Range("I91").Value = 1912.8889
$ws.Range("J91").Value = 2360.25
$ws.Range("K91").Value = 1912.8889
$ws.Range("L91").Value = 2360.25
$ws.Range("M91").Value = -508.8888999999999
$ws.Range("N91").Value = -5168.25
$ws.Range("H97").Value = 503.73685
$ws.Range("I97").Value = 348.5
$ws.Range("K97").Value = 348.5
$ws.Range("M97").Value = 147.5
$ws.Range("H122").Value = 4403.057
$ws.Range("I122").Value = 4634.478
$ws.Range("K122").Value = 13903.434
$ws.Range("M122").Value = -11453.434

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2008.3334
$ws.Range("I86").Value = 1430.8182
$ws.Range("J86").Value = 2405.375
$ws.Range("K86").Value = 1430.8182
$ws.Range("L86").Value = 2405.375
$ws.Range("M86").Value = -307.8181999999999
$ws.Range("N86").Value = -4651.375
$ws.Range("H89").Value = 2008.3334
$ws.Range("I89").Value = 1430.8182
$ws.Range("J89").Value = 2405.375
$ws.Range("K89").Value = 7154.090999999999
$ws.Range("L89").Value = 12026.875
$ws.Range("M89").Value = -1538.090999999999
$ws.Range("N89").Value = -23258.875
$ws.Range("H140").Value = 89999.91
$ws.Range("J140").Value = 89999.91
$ws.Range("L140").Value = 89999.91
$ws.Range("N140").Value = -100359.91

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H58").Value = 2268.2307
$ws.Range("I58").Value = 2471.111
$ws.Range("K58").Value = 2471.111
$ws.Range("M58").Value = -2268.111
$ws.Range("H107").Value = 1180
$ws.Range("I107").Value = 1011
$ws.Range("K107").Value = 1011
$ws.Range("M107").Value = 909
$ws.Range("H136").Value = 2268.2307
$ws.Range("I136").Value = 2471.111
$ws.Range("K136").Value = 7413.333
$ws.Range("M136").Value = -4863.333
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 333333.34
$ws.Range("J141").Value = 333333.34
$ws.Range("L141").Value = 333333.34
$ws.Range("N141").Value = -343693.34

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 207.95653
$ws.Range("I12").Value = 287.45456
$ws.Range("J12").Value = 135.08333
$ws.Range("K12").Value = 862.36368
$ws.Range("L12").Value = 405.24999
$ws.Range("M12").Value = -689.36368
$ws.Range("N12").Value = -751.24999
$ws.Range("H121").Value = 3358.4443
$ws.Range("J121").Value = 4198.5713
$ws.Range("L121").Value = 12595.7139
$ws.Range("N121").Value = -15215.7139
$ws.Range("H137").Value = 2613.2144
$ws.Range("I137").Value = 2343.3333
$ws.Range("J137").Value = 2686.818
$ws.Range("K137").Value = 7029.999899999999
$ws.Range("L137").Value = 8060.454000000001
$ws.Range("M137").Value = -1929.999899999999
$ws.Range("N137").Value = -18260.454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 23334.666
$ws.Range("J20").Value = 24999.5
$ws.Range("L20").Value = 24999.5
$ws.Range("N20").Value = -25489.5
$ws.Range("H126").Value = 4287
$ws.Range("I126").Value = 5313.6
$ws.Range("J126").Value = 3773.7
$ws.Range("K126").Value = 15940.8
$ws.Range("L126").Value = 11321.1
$ws.Range("M126").Value = -13470.8
$ws.Range("N126").Value = -16261.1
$ws.Range("H140").Value = 84812.5
$ws.Range("J140").Value = 84812.5
$ws.Range("L140").Value = 84812.5
$ws.Range("N140").Value = -95172.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2637.25
$ws.Range("I46").Value = 2250
$ws.Range("J46").Value = 3024.5
$ws.Range("K46").Value = 2250
$ws.Range("L46").Value = 3024.5
$ws.Range("M46").Value = -2062
$ws.Range("N46").Value = -3400.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 66665.664
$ws.Range("J34").Value = 99999
$ws.Range("L34").Value = 99999
$ws.Range("N34").Value = -100405
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H43").Value = 32699.5
$ws.Range("I43").Value = 32699.5
$ws.Range("K43").Value = 32699.5
$ws.Range("M43").Value = -32550.5
$ws.Range("H136").Value = 4810.811
$ws.Range("I136").Value = 4283.1562
$ws.Range("K136").Value = 12849.4686
$ws.Range("M136").Value = -10299.4686
